$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2.157254401178704,0.4977145929668723,0.03186391388148557,0.07750110998551385,3.371396017997341,0,0.07973214163530429,0,0.148007748228558,0,0,0.6361323976315703,2.230108892598309),
    @(2.039924398499068,0.4649510196443885,0.03163220842849057,0.0776835693080411,3.329408434567313,0,0.07973214163530429,0,0.1480519533643339,0,0,0.6121101459420117,2.242756627550889),
    @(1.969273045247348,0.4451547636293753,0.03148807221456096,0.07781798545161678,3.305527948551259,0,0.07973214163530429,0,0.1481383786385919,0,0,0.5977641422962847,2.25121920928099),
    @(1.940828943214512,0.4371674501127529,0.03142884699625004,0.07787840105262323,3.296271856290588,0,0.07973214163530429,0,0.148188485101592,0,0,0.592019056531079,2.254842311893121),
    @(1.936126699630904,0.4358459600359481,0.03141898264398169,0.07788877398076011,3.294763533110341,0,0.07973214163530429,0,0.148197703844108,0,0,0.5910711768302335,2.255454446561281),
    @(1.968888036908197,0.4450467218090921,0.03148727548541785,0.07781877738677601,3.305401195893225,0,0.07973214163530429,0,0.1481389941381721,0,0,0.5976862536555387,2.251267365936911),
    @(2.116508977308001,0.4863506388049359,0.03178439813018219,0.07755938231198023,3.356522875535745,0,0.07973214163530429,0,0.1480106735607123,0,0,0.6277654121183645,2.234324789872318),
    @(2.417154895888416,0.569934443296404,0.03235310807379044,0.07722788679950199,3.471966305026484,0,0.07973214163530429,0,0.1482304800106391,0,0,0.6899810506752857,2.206660632456348),
    @(2.645058048405133,0.6329864294455092,0.03276369182069594,0.07709181703288692,3.566223261354537,0,0.07973214163530429,0,0.1486810823754752,0,0,0.7377037233239747,2.189768822909571),
    @(2.7503098524341,0.6620420758732166,0.03294918205976138,0.07705314888685599,3.611193115692458,0,0.07973214163530429,0,0.148949248341296,0,0,0.7598610339173035,2.182839306108662),
    @(2.790396332125795,0.6730994220945377,0.03301926012254697,0.07704183779621232,3.628525898786904,0,0.07973214163530429,0,0.1490599119669724,0,0,0.7683165694784577,2.18032460823342),
    @(2.781752716864162,0.6707155805866023,0.03300417445610471,0.07704412580211439,3.624779421606206,0,0.07973214163530429,0,0.1490356727092887,0,0,0.7664926158031307,2.180861316061169),
    @(2.753603165805316,0.6629506704393293,0.03295495056466713,0.07705215159424661,3.612612991494757,0,0.07973214163530429,0,0.1489581698558951,0,0,0.7605553687665321,2.182630223804338),
    @(2.736390785550043,0.658201578878618,0.03292477893173995,0.07705750125069954,3.605200335126085,0,0.07973214163530429,0,0.1489118850425939,0,0,0.7569271226253846,2.183727997933005),
    @(2.638211599187287,0.6310951623412961,0.03275154557522875,0.07709481078580183,3.563326713419087,0,0.07973214163530429,0,0.148664831003444,0,0,0.7362647588057598,2.190236940317405),
    @(2.578388043800032,0.6145624664045499,0.03264495850221927,0.07712364253578308,3.538176407387084,0,0.07973214163530429,0,0.1485294748784938,0,0,0.7237042523715758,2.194423844988449),
    @(2.54412741306561,0.6050883445625459,0.03258353001874426,0.07714241289024493,3.523907340955191,0,0.07973214163530429,0,0.1484575670595802,0,0,0.7165219138567949,2.19690302216965),
    @(2.532552730047769,0.6018865599619403,0.03256270970499209,0.07714914411893581,3.519109781418734,0,0.07973214163530429,0,0.1484342405713193,0,0,0.7140973152155539,2.197754594449364),
    @(2.584740993375647,0.6163187658634115,0.0326563173900567,0.07712034705577331,3.540833319601859,0,0.07973214163530429,0,0.1485432682305543,0,0,0.7250369747451373,2.193970789588349),
    @(2.761865106983407,0.6652299230225367,0.03296941306427392,0.07704970387233523,3.61617830218924,0,0.07973214163530429,0,0.1489806867075387,0,0,0.762297511896449,2.182107677955344),
    @(2.878967492596018,0.6975149274951491,0.03317309921769507,0.07702295022555106,3.667192022681036,0,0.07973214163530429,0,0.1493197083371243,0,0,0.7870288449404654,2.174992280748782),
    @(2.816343964809107,0.6802543229734397,0.03306446673294339,0.07703545546352863,3.639801994878439,0,0.07973214163530429,0,0.1491338936189592,0,0,0.7737943420799382,2.178731259011528),
    @(2.581868411684297,0.6155246471868736,0.03265118250625321,0.0771218301059573,3.53963153743868,0,0.07973214163530429,0,0.1485370138493067,0,0,0.7244343299538798,2.194175391516495),
    @(2.3346053777978,0.5470391175136342,0.03220062190002704,0.07729865296213134,3.439090955950263,0,0.07973214163530429,0,0.1481203722283979,0,0,0.6727997894977662,2.213544820171961)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

Write-Output "Updated range B2:N25 with new values"